$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A18").Value = "bearbeiten"
$ws.Range("B18").Value = "das Veraendern von existierenden Daten, Objekten,  Eigenschaften oder Verbindungen"

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A18"))
$ws.Sort.SetRange($ws.Range("A2:B18"))
$ws.Sort.Header = 0
$ws.Sort.Apply()
